$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear old rows 2-3 (they get replaced by new content below)
$ws.Range("A2:B3").Clear()

# Row 4: new header-like row
$ws.Range("A4").Value = "%u16"
$ws.Range("B4").Value = "%s"

# Row 5: new data row
$ws.Range("A5").Value = 0
$ws.Range("B5").Value = "<noIssue>"

# Row 6: previously row 2 content, now shifted down
$ws.Range("A6").Value = 1
$ws.Range("B6").Value = "Value out of calibration range"

# Row 7: previously row 3 content, now shifted down
$ws.Range("A7").Value = 2
$ws.Range("B7").Value = "Value below the detection limit (< 1µM)"

$ws.Range("B5").Select()
